$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample row data for rows 7-10 (values only; D/E keep their
# hyperlink/date formatting, same as row 11 already had).
$ws.Range("A7:E10").ClearContents()

# Update the SearchLine selection to D9
$ws.Range("D9").Select()
